$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column values look numeric to Excel (e.g. "1.00", "0.0210") and would
# otherwise be auto-converted to actual numbers, losing their original text formatting
# (trailing zeros, leading zeros, scientific notation, etc). Force those particular
# cells to text via NumberFormat before assigning, then restore the default "Normal"
# cell style so no stray formatting is left behind.

$ws.Range("D2").Value = '37.449.69'
$ws.Range("E2").Value = '  +1.19%  '

$ws.Range("D3").Value = '2.077.01'
$ws.Range("E3").Value = '  +4.27%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.91%  '

$ws.Range("E6").Value = '  +2.16%  '

$ws.Range("B7").Value = 'Solana'
$ws.Range("C7").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.49'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.86%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +2.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.25%  '

$ws.Range("E11").Value = '  +0.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.102'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.77%  '

$ws.Range("D13").Value = '2.382.51'
$ws.Range("E13").Value = '  +4.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.775'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.53%  '

$ws.Range("E17").Value = '  +4.20%  '

$ws.Range("D18").Value = '2.109.17'
$ws.Range("E18").Value = '  +5.85%  '

$ws.Range("D19").Value = '37.570.40'
$ws.Range("E19").Value = '  +1.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +20.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").Value = '0.0₃0811'
$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '223.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.97%  '

$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("E25").Value = '  +4.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.77%  '

$ws.Range("E28").Value = '  +2.55%  '

$ws.Range("E29").Value = '  +6.75%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.118'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.46'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0623'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.40'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.37%  '

$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +13.92%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.86%  '

$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.24%  '

$ws.Range("E41").Value = '  -4.72%  '

$ws.Range("D43").Value = '1.476.24'
$ws.Range("E43").Value = '  +2.95%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +18.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '94.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.57%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.27%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0210'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.34%  '

$ws.Range("E48").Value = '  +1.27%  '

$ws.Range("E49").Value = '  +2.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.17%  '

$ws.Range("E51").Value = '  +1.65%  '
